$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Replace the Latitude/Longitude columns (N/O) with a single "Maps" link
#    column. N keeps its letter; the new link text goes into N. The old
#    "Longitude" column (O) together with the following blank placeholder
#    column (P) are removed outright, which shifts the old Q:U placeholder
#    columns left into O:S (they are already blank, so nothing else to do
#    there).
# ---------------------------------------------------------------------------

$mapsUrl = @{
    2  = "https://goo.gl/maps/d4YmuG9wxEysZxWJ8"
    3  = "https://goo.gl/maps/9Ts4xERRLZsuQzVE8"
    4  = "https://goo.gl/maps/mH8Z3SF8rH3Ynj7MA"
    5  = "https://goo.gl/maps/poLHvL9UBWTKNb1G8"
    6  = "https://goo.gl/maps/knUTS44tkG36Ppdt7"
    7  = "https://goo.gl/maps/tyeThCtvT7YVXXux8"
    8  = "https://goo.gl/maps/VHvmVFozaqfZXvbE9"
    9  = "https://goo.gl/maps/zWZDt7W5PCc5xtP89"
    10 = "https://goo.gl/maps/tWDMyawW21QqTiaU6"
    11 = "https://goo.gl/maps/YkL7HMAym5shWYoAA"
    12 = "https://goo.gl/maps/pgfRRWX9ei9EKHHdA"
    13 = "https://goo.gl/maps/uB7VXqb4VMsjR74n8"
    14 = "https://goo.gl/maps/M5vbmKxtnagZ3Jp66"
    15 = "https://goo.gl/maps/S4GX5oriaQ5forUB9"
    16 = "https://goo.gl/maps/D4Ptn8UrNXR4Bbkj8"
    17 = "https://goo.gl/maps/rgv2MC5ffMEy2F4L6"
    18 = "https://goo.gl/maps/zWZDt7W5PCc5xtP89"
    19 = "https://goo.gl/maps/yHzn96oNnfMyrZP4A"
    20 = "https://goo.gl/maps/NctRQb1oQU2F7XcC9"
    21 = "https://goo.gl/maps/r9G6Uno2YAgA4LEg7"
    22 = "https://goo.gl/maps/nSg9FFxLRmcBg3Z18"
    23 = "https://goo.gl/maps/eUHfiCPwgfRi7gae6"
    24 = "https://goo.gl/maps/fCPpAtTWQ5icqfT66"
    25 = "https://goo.gl/maps/VkjREyA6zNJcW6dL6"
    26 = "https://goo.gl/maps/43LkWYeDxApzgXTH8"
    27 = "https://goo.gl/maps/TVoeSWVHrYohc9Ra6"
}

# Header
$ws.Range("N1").Value = "Maps"

# Data rows 2..27 (row 8 gets a live hyperlink further below)
foreach ($r in $mapsUrl.Keys) {
    $ws.Range("N" + $r).Value = $mapsUrl[$r]
}

# Normalise every cell in column N (header + data + the couple of blank
# rows below the table) back to the plain "Normal" look instead of the old
# text-number format - mirror the look already used by column M.
$ws.Range("M1:M30").Copy()
$ws.Columns("N").PasteSpecial(-4122)

# Turn N8 into a real clickable hyperlink (this also is the one cell that
# keeps the underlined "Hyperlink" look).
$ws.Hyperlinks.Add($ws.Range("N8"), $mapsUrl[8])

# N5 visually uses the same "Hyperlink" style as N8 (copy the look only,
# not an actual link).
$ws.Range("N8").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = $mapsUrl[5]

# Drop the old Longitude column (O) and the blank column after it (P);
# everything to the right shifts left and is already blank so it lines up
# with the new O:S placeholder columns.
$ws.Columns("O:P").Delete()

# ---------------------------------------------------------------------------
# 2. Shrink the filter / used range from column U to column S everywhere.
# ---------------------------------------------------------------------------
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$S`$27"

$ws.AutoFilterMode = $false
$ws.Range("A1:S27").AutoFilter()

# ---------------------------------------------------------------------------
# 3. Misc housekeeping to match the saved state.
# ---------------------------------------------------------------------------
$ws.Range("N30").Select()

Write-Output "done"
